$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("G3").Value = 2.4
$ws.Range("H3").Value = 3.15
$ws.Range("I3").Value = 2.85
$ws.Range("L3").Value = 1.44
$ws.Range("M3").Value = 2.42
$ws.Range("N3").Value = 2.25
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 1.5
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.98
$ws.Range("S3").Value = 1.65
$ws.Range("T3").Value = 6.4
$ws.Range("U3").Value = 10.25
$ws.Range("V3").Value = 10
$ws.Range("W3").Value = 24
$ws.Range("X3").Value = 23
$ws.Range("Y3").Value = 40
$ws.Range("Z3").Value = 7.3
$ws.Range("AA3").Value = 6.2
$ws.Range("AB3").Value = 18
$ws.Range("AC3").Value = 110
$ws.Range("AE3").Value = 7.2
$ws.Range("AF3").Value = 13
$ws.Range("AG3").Value = 11
$ws.Range("AH3").Value = 32
$ws.Range("AI3").Value = 29
$ws.Range("AJ3").Value = 45

# Row 4
$ws.Range("G4").Value = 1.65
$ws.Range("H4").Value = 3.45
$ws.Range("I4").Value = 5.1
$ws.Range("L4").Value = 1.42
$ws.Range("M4").Value = 2.47
$ws.Range("N4").Value = 2.2
$ws.Range("O4").Value = 1.52
$ws.Range("P4").Value = 1.5
$ws.Range("Q4").Value = 2.27
$ws.Range("R4").Value = 2.18
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 5.1
$ws.Range("U4").Value = 6.5
$ws.Range("V4").Value = 8.75
$ws.Range("W4").Value = 11.75
$ws.Range("X4").Value = 16
$ws.Range("Y4").Value = 40
$ws.Range("Z4").Value = 7.5
$ws.Range("AA4").Value = 7
$ws.Range("AB4").Value = 23
$ws.Range("AC4").Value = 150
$ws.Range("AE4").Value = 10.5
$ws.Range("AF4").Value = 26
$ws.Range("AG4").Value = 18
$ws.Range("AH4").Value = 100
$ws.Range("AI4").Value = 65
$ws.Range("AJ4").Value = 80

# Row 5
$ws.Range("G5").Value = 1.98
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 3.8
$ws.Range("L5").Value = 1.4
$ws.Range("M5").Value = 2.5
$ws.Range("N5").Value = 2.18
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 1.45
$ws.Range("Q5").Value = 2.37
$ws.Range("R5").Value = 1.93
$ws.Range("S5").Value = 1.7
$ws.Range("T5").Value = 5.9
$ws.Range("U5").Value = 8.5
$ws.Range("V5").Value = 8.75
$ws.Range("W5").Value = 17.5
$ws.Range("X5").Value = 18
$ws.Range("Y5").Value = 35
$ws.Range("Z5").Value = 7.4
$ws.Range("AA5").Value = 6.1
$ws.Range("AB5").Value = 17
$ws.Range("AC5").Value = 100
$ws.Range("AE5").Value = 9
$ws.Range("AF5").Value = 19.5
$ws.Range("AG5").Value = 13.5
$ws.Range("AH5").Value = 60
$ws.Range("AI5").Value = 40
$ws.Range("AJ5").Value = 55

# Row 6
$ws.Range("G6").Value = 2.55
$ws.Range("H6").Value = 3.25
$ws.Range("I6").Value = 2.57
$ws.Range("L6").Value = 1.37
$ws.Range("M6").Value = 2.62
$ws.Range("N6").Value = 2.07
$ws.Range("O6").Value = 1.6
$ws.Range("P6").Value = 1.42
$ws.Range("Q6").Value = 2.45
$ws.Range("R6").Value = 1.85
$ws.Range("S6").Value = 1.75
$ws.Range("T6").Value = 7.3
$ws.Range("U6").Value = 11.75
$ws.Range("V6").Value = 10
$ws.Range("W6").Value = 27
$ws.Range("X6").Value = 23
$ws.Range("Y6").Value = 37
$ws.Range("Z6").Value = 8.25
$ws.Range("AA6").Value = 6.3
$ws.Range("AB6").Value = 16
$ws.Range("AC6").Value = 90
$ws.Range("AD6").Value = 800
$ws.Range("AE6").Value = 7.4
$ws.Range("AF6").Value = 11.75
$ws.Range("AG6").Value = 10
$ws.Range("AH6").Value = 28
$ws.Range("AI6").Value = 24
$ws.Range("AJ6").Value = 37

